# Add new row of data (row 97) to Sheet1, as committed via Streamlit on 2024-12-04 11:21:19
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 97

$ws.Cells.Item($row, 1).Value = "Kindergarden"
$ws.Cells.Item($row, 2).Value = "Kindergarden Almere Zwitserlandstraat"
$ws.Cells.Item($row, 3).Value = "KDV"

# Keep the report date as plain text (matches the rest of column D), not an
# auto-converted date serial number.
$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "2023-11-21"
$ws.Cells.Item($row, 4).Style = "Normal"

$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 1
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
